$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("B18").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("B23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("H31").Value = 0.09399999999999997
$ws.Range("E32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K43").Value = -0.111
$ws.Range("C44").Value = -0.02100000000000002
$ws.Range("I44").Value = 0.009000000000000008
$ws.Range("G46").Value = 0
$ws.Range("G48").Value = 0.007000000000000006
$ws.Range("K48").Value = -0.02699999999999997
$ws.Range("C49").Value = 0
$ws.Range("E56").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("E60").Value = -0.02299999999999991
$ws.Range("J61").Value = 0
$ws.Range("G62").Value = 0.103
$ws.Range("C65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("H68").Value = 0.062
$ws.Range("I68").Value = 0.07299999999999995
$ws.Range("G69").Value = 0.02900000000000003
$ws.Range("K69").Value = -0.02899999999999997
$ws.Range("K71").Value = -0.04699999999999999
$ws.Range("E73").Value = -0.04100000000000004
$ws.Range("D76").Value = 0.02100000000000002
$ws.Range("G77").Value = 0.04299999999999993
$ws.Range("F79").Value = -0.01300000000000001
$ws.Range("L79").Value = -0.02000000000000002
$ws.Range("I80").Value = 0.04899999999999993
$ws.Range("G81").Value = 0.09599999999999997
$ws.Range("I81").Value = 0.03400000000000003
$ws.Range("G82").Value = 0
$ws.Range("J84").Value = 0.147
$ws.Range("I86").Value = 0.04500000000000004
$ws.Range("K87").Value = -0.122
$ws.Range("L87").Value = -0.03499999999999998
$ws.Range("G88").Value = 0.05100000000000005
$ws.Range("I89").Value = 0.127
$ws.Range("C90").Value = -0.02000000000000002
$ws.Range("K90").Value = -0.08700000000000002
$ws.Range("C91").Value = -0.02900000000000003
$ws.Range("E91").Value = -0.03300000000000003
$ws.Range("J93").Value = 0.111
$ws.Range("G95").Value = 0.04900000000000004
$ws.Range("C97").Value = -0.004000000000000004
$ws.Range("G97").Value = 0.08000000000000007
$ws.Range("H98").Value = 0.03300000000000003
$ws.Range("D99").Value = 0.03199999999999997
$ws.Range("M99").Value = 2.758823529411765
$ws.Range("G100").Value = 0.06100000000000005
